$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: the inline-data PNG placeholder is replaced by the real Soriana image URL.
$ws.Range("D14").Value = "https://www.soriana.com/dw/image/v2/BGBD_PRD/on/demandware.static/-/Sites-soriana-grocery-master-catalog/default/dwaaae9371/images/product/8806091641540_A.jpg?sw=106&sh=106&sm=fit"

# The "Hisense Smart TV LED A6H 55"" row (old row 94) was removed from the
# product list; every following TV row shifts up by one, and a brand-new
# "Xiaomi Smart TV LED A Pro 55 55""" row is appended at the end (row 109).
$ws.Rows("94").Delete()

$ws.Range("A109").Value = "Xiaomi Smart TV LED A Pro 55 55"", 4K Ultra HD, Negro"
$ws.Range("B109").Value = 7689
$ws.Range("C109").Value = "3840 x 2160"
$ws.Range("D109").Value = "https://www.cyberpuerta.mx/img/product/S/CP-XIAOMI-L55M8-A2LA-0b4b8c.jpg"
$ws.Range("E109").Value = "Cyberpuerta"
$ws.Range("F109").Value = ": 139,7 cm (55"")"
